$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Typo / punctuation fixes (simple Find & Replace, scoped to the
#    whole document content; each search string is unique).
# ------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("Quy trình hoạt động của nhà hang", $true, $false, $false, $false, $false, $true, 1, $false, "Quy trình hoạt động của nhà hàng", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Sau khí khách hàng đến", $true, $false, $false, $false, $false, $true, 1, $false, "Sau khi khách hàng đến", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("khách có thể phục vụ", $true, $false, $false, $false, $false, $true, 1, $false, "khách có thể phục vụ.", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("khách hàng có thể đặt món từ trước", $true, $false, $false, $false, $false, $true, 1, $false, "khách hàng có thể đặt món từ trước.", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Hiển thị sơ đồ nhà hàng vơi các trạng thái", $true, $false, $false, $false, $false, $true, 1, $false, "Hiển thị sơ đồ nhà hàng với  các trạng thái", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Remove the whole paragraph "Xác định vị trí bàn."
# ------------------------------------------------------------------

$rng = $d.Content
$found = $rng.Find.Execute("Xác định vị trí bàn.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Expand(4)
    $rng.Delete()
}

# ------------------------------------------------------------------
# 3) Remove the (second) empty "ListParagraph" paragraph that sits
#    between "...các món ăn." and the ind=1440 empty paragraph, right
#    before "Các chức năng cần thiết của ứng dụng quản lý nhà hàng".
# ------------------------------------------------------------------

$rng = $d.Content
$found = $rng.Find.Execute("Nhà hàng cần quản lý kho nguyên liệu để đảm bảo luôn có đủ nguyên liệu cần thiết để chế biến các món ăn.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.Move(4, 2)
    $rng.Expand(4)
    $rng.Delete()
}

# ------------------------------------------------------------------
# 4) Re-style the blank paragraph that follows the title, then insert
#    a new "Nhóm 5 - Ứng dụng quản lý nhà hàng" paragraph after it.
# ------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$blankPara = $titlePara.Next()
$blankPara.Alignment = 1
$startPos = $blankPara.Range.Start
$blankPara.Range.Text = "X"
$blankPara.Range.Font.Size = 16
$blankPara.Range.Font.SizeBi = 16
$d.Range($startPos, $startPos + 1).Delete()

$blankPara = $titlePara.Next()
$blankPara.Range.InsertParagraphAfter()
$newPara = $blankPara.Next()
$newPara.Range.Text = "Nhóm 5 - Ứng dụng quản lý nhà hàng"
